# Simplify language by using precedence rules
# Appends a trailing semicolon to the various Action/Query expressions in the
# "ExampleProcess" workbook (letting precedence rules make the code implicit),
# and converts the rich-text "GO(...)" cells into plain text while they're at it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Action column (E) updates -------------------------------------------------
$ws.Range("E3").Value = "GO(PreviousMedBrand);"
$ws.Range("E4").Value = "GO(PreviousMedBrand);"
$ws.Range("E5").Value = "GO(PreviousMedBrand);"
$ws.Range("E6").Value = "GO(PreviousMedProduct);"
$ws.Range("E7").Value = "GO(SeeProduct);"
$ws.Range("E8").Value = "GO(FreeText);"
$ws.Range("E9").Value = "FINISH();"

# --- Query column (C) updates ---------------------------------------------------
$ws.Range("C6").Value = "SELECT * FROM Brands;"
$ws.Range("C7").Value = "SELECT * FROM Products WHERE brand_id == [prev_med_brand];"
$ws.Range("C8").Value = "SELECT * FROM Products WHERE id == [med];"

# --- Row height tweaks for rows 7 & 8 -------------------------------------------
$ws.Rows.Item(7).RowHeight = 13.8
$ws.Rows.Item(8).RowHeight = 13.8

# --- Restore view to top-left and move the active selection to C14 -------------
[void]$ws.Application.Goto($ws.Range("A1"))
[void]$ws.Range("C14").Select()
